$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP1912")

# Header updates
$ws.Cells.Item(2, 1).Value = 'Última actualización: 16:43:37'
$ws.Cells.Item(3, 1).Value = 'Total filas: 379'

# Row data updates (rows 39-384), reflecting refreshed scrape/sort results
# Row 39
$ws.Cells.Item(39, 1).Value = '06:43:40'
$ws.Cells.Item(39, 2).Value = '06:46'
$ws.Cells.Item(39, 3).Value = '225_C ROCA-H SUR'
$ws.Cells.Item(39, 4).Value = 3
$ws.Cells.Item(39, 5).Value = 'LP1912'

# Row 40
$ws.Cells.Item(40, 1).Value = '05:18:56'
$ws.Cells.Item(40, 2).Value = '06:46'
$ws.Cells.Item(40, 3).Value = '215C_EL PATO'
$ws.Cells.Item(40, 4).Value = 88
$ws.Cells.Item(40, 5).Value = 'LP1912'

# Row 111
$ws.Cells.Item(111, 1).Value = '07:47:32'
$ws.Cells.Item(111, 2).Value = '09:23'
$ws.Cells.Item(111, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(111, 4).Value = 96
$ws.Cells.Item(111, 5).Value = 'LP1912'

# Row 112
$ws.Cells.Item(112, 1).Value = '08:57:13'
$ws.Cells.Item(112, 2).Value = '09:23'
$ws.Cells.Item(112, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(112, 4).Value = 26
$ws.Cells.Item(112, 5).Value = 'LP1912'

# Row 121
$ws.Cells.Item(121, 1).Value = '09:38:09'
$ws.Cells.Item(121, 2).Value = '09:41'
$ws.Cells.Item(121, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(121, 4).Value = 3
$ws.Cells.Item(121, 5).Value = 'LP1912'

# Row 122
$ws.Cells.Item(122, 1).Value = '08:21:50'
$ws.Cells.Item(122, 2).Value = '09:41'
$ws.Cells.Item(122, 3).Value = '215C_EL PATO'
$ws.Cells.Item(122, 4).Value = 80
$ws.Cells.Item(122, 5).Value = 'LP1912'

# Row 123
$ws.Cells.Item(123, 1).Value = '09:38:09'
$ws.Cells.Item(123, 2).Value = '09:41'
$ws.Cells.Item(123, 3).Value = '14_ABASTO'
$ws.Cells.Item(123, 4).Value = 3
$ws.Cells.Item(123, 5).Value = 'LP1912'

# Row 206
$ws.Cells.Item(206, 1).Value = '11:48:04'
$ws.Cells.Item(206, 2).Value = '12:35'
$ws.Cells.Item(206, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(206, 4).Value = 47
$ws.Cells.Item(206, 5).Value = 'LP1912'

# Row 207
$ws.Cells.Item(207, 1).Value = '11:48:04'
$ws.Cells.Item(207, 2).Value = '12:35'
$ws.Cells.Item(207, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(207, 4).Value = 47
$ws.Cells.Item(207, 5).Value = 'LP1912'

# Row 231
$ws.Cells.Item(231, 1).Value = '11:48:04'
$ws.Cells.Item(231, 2).Value = '13:21'
$ws.Cells.Item(231, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(231, 4).Value = 93
$ws.Cells.Item(231, 5).Value = 'LP1912'

# Row 232
$ws.Cells.Item(232, 1).Value = '13:19:56'
$ws.Cells.Item(232, 2).Value = '13:21'
$ws.Cells.Item(232, 3).Value = '10_OLMOS'
$ws.Cells.Item(232, 4).Value = 2
$ws.Cells.Item(232, 5).Value = 'LP1912'

# Row 242
$ws.Cells.Item(242, 1).Value = '13:19:56'
$ws.Cells.Item(242, 2).Value = '13:46'
$ws.Cells.Item(242, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(242, 4).Value = 27
$ws.Cells.Item(242, 5).Value = 'LP1912'

# Row 243
$ws.Cells.Item(243, 1).Value = '11:48:04'
$ws.Cells.Item(243, 2).Value = '13:46'
$ws.Cells.Item(243, 3).Value = '17_ROMERO'
$ws.Cells.Item(243, 4).Value = 118
$ws.Cells.Item(243, 5).Value = 'LP1912'

# Row 253
$ws.Cells.Item(253, 1).Value = '13:53:08'
$ws.Cells.Item(253, 2).Value = '13:57'
$ws.Cells.Item(253, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(253, 4).Value = 4
$ws.Cells.Item(253, 5).Value = 'LP1912'

# Row 254
$ws.Cells.Item(254, 1).Value = '12:37:14'
$ws.Cells.Item(254, 2).Value = '13:57'
$ws.Cells.Item(254, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(254, 4).Value = 80
$ws.Cells.Item(254, 5).Value = 'LP1912'

# Row 278
$ws.Cells.Item(278, 1).Value = '14:58:43'
$ws.Cells.Item(278, 2).Value = '14:58'
$ws.Cells.Item(278, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(278, 4).Value = 0
$ws.Cells.Item(278, 5).Value = 'LP1912'

# Row 279
$ws.Cells.Item(279, 1).Value = '13:19:56'
$ws.Cells.Item(279, 2).Value = '14:58'
$ws.Cells.Item(279, 3).Value = '215B_EL PATO'
$ws.Cells.Item(279, 4).Value = 99
$ws.Cells.Item(279, 5).Value = 'LP1912'

# Row 312
$ws.Cells.Item(312, 1).Value = '15:57:48'
$ws.Cells.Item(312, 2).Value = '15:57'
$ws.Cells.Item(312, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(312, 4).Value = 0
$ws.Cells.Item(312, 5).Value = 'LP1912'

# Row 313
$ws.Cells.Item(313, 1).Value = '15:57:48'
$ws.Cells.Item(313, 2).Value = '15:57'
$ws.Cells.Item(313, 3).Value = '17_ROMERO'
$ws.Cells.Item(313, 4).Value = 0
$ws.Cells.Item(313, 5).Value = 'LP1912'

# Row 323
$ws.Cells.Item(323, 1).Value = '15:31:33'
$ws.Cells.Item(323, 2).Value = '16:20'
$ws.Cells.Item(323, 3).Value = '215C_EL PATO'
$ws.Cells.Item(323, 4).Value = 49
$ws.Cells.Item(323, 5).Value = 'LP1912'

# Row 324
$ws.Cells.Item(324, 1).Value = '14:46:52'
$ws.Cells.Item(324, 2).Value = '16:20'
$ws.Cells.Item(324, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(324, 4).Value = 94
$ws.Cells.Item(324, 5).Value = 'LP1912'

# Row 343
$ws.Cells.Item(343, 1).Value = '16:43:37'
$ws.Cells.Item(343, 2).Value = '16:44'
$ws.Cells.Item(343, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(343, 4).Value = 1
$ws.Cells.Item(343, 5).Value = 'LP1912'

# Row 344
$ws.Cells.Item(344, 1).Value = '15:31:33'
$ws.Cells.Item(344, 2).Value = '16:48'
$ws.Cells.Item(344, 3).Value = '15_ABASTO'
$ws.Cells.Item(344, 4).Value = 77
$ws.Cells.Item(344, 5).Value = 'LP1912'

# Row 345
$ws.Cells.Item(345, 1).Value = '15:57:48'
$ws.Cells.Item(345, 2).Value = '16:50'
$ws.Cells.Item(345, 3).Value = '14_ABASTO'
$ws.Cells.Item(345, 4).Value = 53
$ws.Cells.Item(345, 5).Value = 'LP1912'

# Row 346
$ws.Cells.Item(346, 1).Value = '16:33:08'
$ws.Cells.Item(346, 2).Value = '16:51'
$ws.Cells.Item(346, 3).Value = '14_ABASTO'
$ws.Cells.Item(346, 4).Value = 18
$ws.Cells.Item(346, 5).Value = 'LP1912'

# Row 347
$ws.Cells.Item(347, 1).Value = '14:58:43'
$ws.Cells.Item(347, 2).Value = '16:56'
$ws.Cells.Item(347, 3).Value = '17_179 Y 38'
$ws.Cells.Item(347, 4).Value = 118
$ws.Cells.Item(347, 5).Value = 'LP1912'

# Row 348
$ws.Cells.Item(348, 1).Value = '16:18:55'
$ws.Cells.Item(348, 2).Value = '16:56'
$ws.Cells.Item(348, 3).Value = '10_OLMOS'
$ws.Cells.Item(348, 4).Value = 38
$ws.Cells.Item(348, 5).Value = 'LP1912'

# Row 349
$ws.Cells.Item(349, 1).Value = '16:33:08'
$ws.Cells.Item(349, 2).Value = '16:57'
$ws.Cells.Item(349, 3).Value = '10_OLMOS'
$ws.Cells.Item(349, 4).Value = 24
$ws.Cells.Item(349, 5).Value = 'LP1912'

# Row 350
$ws.Cells.Item(350, 1).Value = '15:57:48'
$ws.Cells.Item(350, 2).Value = '17:04'
$ws.Cells.Item(350, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(350, 4).Value = 67
$ws.Cells.Item(350, 5).Value = 'LP1912'

# Row 351
$ws.Cells.Item(351, 1).Value = '16:18:55'
$ws.Cells.Item(351, 2).Value = '17:04'
$ws.Cells.Item(351, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(351, 4).Value = 46
$ws.Cells.Item(351, 5).Value = 'LP1912'

# Row 352
$ws.Cells.Item(352, 1).Value = '15:31:33'
$ws.Cells.Item(352, 2).Value = '17:04'
$ws.Cells.Item(352, 3).Value = '215A_EL PATO'
$ws.Cells.Item(352, 4).Value = 93
$ws.Cells.Item(352, 5).Value = 'LP1912'

# Row 353
$ws.Cells.Item(353, 1).Value = '16:33:08'
$ws.Cells.Item(353, 2).Value = '17:05'
$ws.Cells.Item(353, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(353, 4).Value = 32
$ws.Cells.Item(353, 5).Value = 'LP1912'

# Row 355
$ws.Cells.Item(355, 1).Value = '16:33:08'
$ws.Cells.Item(355, 2).Value = '17:05'
$ws.Cells.Item(355, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(355, 4).Value = 32
$ws.Cells.Item(355, 5).Value = 'LP1912'

# Row 356
$ws.Cells.Item(356, 1).Value = '16:33:08'
$ws.Cells.Item(356, 2).Value = '17:10'
$ws.Cells.Item(356, 3).Value = '10_OLMOS'
$ws.Cells.Item(356, 4).Value = 37
$ws.Cells.Item(356, 5).Value = 'LP1912'

# Row 357
$ws.Cells.Item(357, 1).Value = '16:43:37'
$ws.Cells.Item(357, 2).Value = '17:16'
$ws.Cells.Item(357, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(357, 4).Value = 33
$ws.Cells.Item(357, 5).Value = 'LP1912'

# Row 358
$ws.Cells.Item(358, 1).Value = '16:18:55'
$ws.Cells.Item(358, 2).Value = '17:20'
$ws.Cells.Item(358, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(358, 4).Value = 62
$ws.Cells.Item(358, 5).Value = 'LP1912'

# Row 359
$ws.Cells.Item(359, 1).Value = '16:18:55'
$ws.Cells.Item(359, 2).Value = '17:20'
$ws.Cells.Item(359, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(359, 4).Value = 62
$ws.Cells.Item(359, 5).Value = 'LP1912'

# Row 360
$ws.Cells.Item(360, 1).Value = '15:31:33'
$ws.Cells.Item(360, 2).Value = '17:21'
$ws.Cells.Item(360, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(360, 4).Value = 110
$ws.Cells.Item(360, 5).Value = 'LP1912'

# Row 361
$ws.Cells.Item(361, 1).Value = '16:33:08'
$ws.Cells.Item(361, 2).Value = '17:21'
$ws.Cells.Item(361, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(361, 4).Value = 48
$ws.Cells.Item(361, 5).Value = 'LP1912'

# Row 362
$ws.Cells.Item(362, 1).Value = '15:31:33'
$ws.Cells.Item(362, 2).Value = '17:24'
$ws.Cells.Item(362, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(362, 4).Value = 113
$ws.Cells.Item(362, 5).Value = 'LP1912'

# Row 363
$ws.Cells.Item(363, 1).Value = '16:18:55'
$ws.Cells.Item(363, 2).Value = '17:28'
$ws.Cells.Item(363, 3).Value = '14_ABASTO'
$ws.Cells.Item(363, 4).Value = 70
$ws.Cells.Item(363, 5).Value = 'LP1912'

# Row 364
$ws.Cells.Item(364, 1).Value = '16:33:08'
$ws.Cells.Item(364, 2).Value = '17:29'
$ws.Cells.Item(364, 3).Value = '14_ABASTO'
$ws.Cells.Item(364, 4).Value = 56
$ws.Cells.Item(364, 5).Value = 'LP1912'

# Row 365
$ws.Cells.Item(365, 1).Value = '16:18:55'
$ws.Cells.Item(365, 2).Value = '17:30'
$ws.Cells.Item(365, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(365, 4).Value = 72
$ws.Cells.Item(365, 5).Value = 'LP1912'

# Row 366
$ws.Cells.Item(366, 1).Value = '16:33:08'
$ws.Cells.Item(366, 2).Value = '17:31'
$ws.Cells.Item(366, 3).Value = '15_ABASTO'
$ws.Cells.Item(366, 4).Value = 58
$ws.Cells.Item(366, 5).Value = 'LP1912'

# Row 367
$ws.Cells.Item(367, 1).Value = '16:33:08'
$ws.Cells.Item(367, 2).Value = '17:33'
$ws.Cells.Item(367, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(367, 4).Value = 60
$ws.Cells.Item(367, 5).Value = 'LP1912'

# Row 368
$ws.Cells.Item(368, 1).Value = '16:43:37'
$ws.Cells.Item(368, 2).Value = '17:34'
$ws.Cells.Item(368, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(368, 4).Value = 51
$ws.Cells.Item(368, 5).Value = 'LP1912'

# Row 369
$ws.Cells.Item(369, 1).Value = '15:57:48'
$ws.Cells.Item(369, 2).Value = '17:35'
$ws.Cells.Item(369, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(369, 4).Value = 98
$ws.Cells.Item(369, 5).Value = 'LP1912'

# Row 370
$ws.Cells.Item(370, 1).Value = '16:43:37'
$ws.Cells.Item(370, 2).Value = '17:36'
$ws.Cells.Item(370, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(370, 4).Value = 53
$ws.Cells.Item(370, 5).Value = 'LP1912'

# Row 371
$ws.Cells.Item(371, 1).Value = '16:18:55'
$ws.Cells.Item(371, 2).Value = '17:38'
$ws.Cells.Item(371, 3).Value = '17_ROMERO'
$ws.Cells.Item(371, 4).Value = 80
$ws.Cells.Item(371, 5).Value = 'LP1912'

# Row 372
$ws.Cells.Item(372, 1).Value = '16:18:55'
$ws.Cells.Item(372, 2).Value = '17:39'
$ws.Cells.Item(372, 3).Value = '215B_EL PATO'
$ws.Cells.Item(372, 4).Value = 81
$ws.Cells.Item(372, 5).Value = 'LP1912'

# Row 373
$ws.Cells.Item(373, 1).Value = '15:57:48'
$ws.Cells.Item(373, 2).Value = '17:40'
$ws.Cells.Item(373, 3).Value = '215B_EL PATO'
$ws.Cells.Item(373, 4).Value = 103
$ws.Cells.Item(373, 5).Value = 'LP1912'

# Row 374
$ws.Cells.Item(374, 1).Value = '16:43:37'
$ws.Cells.Item(374, 2).Value = '17:40'
$ws.Cells.Item(374, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(374, 4).Value = 57
$ws.Cells.Item(374, 5).Value = 'LP1912'

# Row 375
$ws.Cells.Item(375, 1).Value = '15:57:48'
$ws.Cells.Item(375, 2).Value = '17:41'
$ws.Cells.Item(375, 3).Value = '17_ROMERO'
$ws.Cells.Item(375, 4).Value = 104
$ws.Cells.Item(375, 5).Value = 'LP1912'

# Row 376
$ws.Cells.Item(376, 1).Value = '16:33:08'
$ws.Cells.Item(376, 2).Value = '17:41'
$ws.Cells.Item(376, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(376, 4).Value = 68
$ws.Cells.Item(376, 5).Value = 'LP1912'

# Row 377
$ws.Cells.Item(377, 1).Value = '16:43:37'
$ws.Cells.Item(377, 2).Value = '17:45'
$ws.Cells.Item(377, 3).Value = '15_ABASTO'
$ws.Cells.Item(377, 4).Value = 62
$ws.Cells.Item(377, 5).Value = 'LP1912'

# Row 378
$ws.Cells.Item(378, 1).Value = '15:57:48'
$ws.Cells.Item(378, 2).Value = '17:50'
$ws.Cells.Item(378, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(378, 4).Value = 113
$ws.Cells.Item(378, 5).Value = 'LP1912'

# Row 379
$ws.Cells.Item(379, 1).Value = '16:33:08'
$ws.Cells.Item(379, 2).Value = '17:51'
$ws.Cells.Item(379, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(379, 4).Value = 78
$ws.Cells.Item(379, 5).Value = 'LP1912'

# Row 380
$ws.Cells.Item(380, 1).Value = '15:57:48'
$ws.Cells.Item(380, 2).Value = '17:52'
$ws.Cells.Item(380, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(380, 4).Value = 115
$ws.Cells.Item(380, 5).Value = 'LP1912'

# Row 381
$ws.Cells.Item(381, 1).Value = '16:33:08'
$ws.Cells.Item(381, 2).Value = '18:04'
$ws.Cells.Item(381, 3).Value = '17_ROMERO'
$ws.Cells.Item(381, 4).Value = 91
$ws.Cells.Item(381, 5).Value = 'LP1912'

# Row 382
$ws.Cells.Item(382, 1).Value = '16:33:08'
$ws.Cells.Item(382, 2).Value = '18:21'
$ws.Cells.Item(382, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(382, 4).Value = 108
$ws.Cells.Item(382, 5).Value = 'LP1912'

# Row 383
$ws.Cells.Item(383, 1).Value = '16:33:08'
$ws.Cells.Item(383, 2).Value = '18:28'
$ws.Cells.Item(383, 3).Value = '215C_EL PATO'
$ws.Cells.Item(383, 4).Value = 115
$ws.Cells.Item(383, 5).Value = 'LP1912'

# Row 384
$ws.Cells.Item(384, 1).Value = '16:43:37'
$ws.Cells.Item(384, 2).Value = '18:32'
$ws.Cells.Item(384, 3).Value = '11X44_ETCHEVERRY'
$ws.Cells.Item(384, 4).Value = 109
$ws.Cells.Item(384, 5).Value = 'LP1912'

# Update "Última actualización" timestamp on the other two sheets
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = 'Última actualización: 16:43:37'
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = 'Última actualización: 16:43:37'
